$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 224
$ws.Range("I9").Value = 100
$ws.Range("K9").Value = 100
$ws.Range("M9").Value = 69

$ws.Range("H39").Value = 272.33334
$ws.Range("I39").Value = 83
$ws.Range("J39").Value = 651
$ws.Range("K39").Value = 249
$ws.Range("L39").Value = 1953
$ws.Range("M39").Value = 47
$ws.Range("N39").Value = -2545

$ws.Range("H137").Value = 2271.4644
$ws.Range("I137").Value = 1954.762
$ws.Range("J137").Value = 3221.5715
$ws.Range("K137").Value = 5864.286
$ws.Range("L137").Value = 9664.7145
$ws.Range("M137").Value = -3314.286
$ws.Range("N137").Value = -14764.7145

$ws.Range("H138").Value = 2159
$ws.Range("J138").Value = 3911.5
$ws.Range("L138").Value = 11734.5
$ws.Range("N138").Value = -22014.5

$ws.Range("H140").Value = 85277.75
$ws.Range("J140").Value = 85277.75
$ws.Range("L140").Value = 85277.75
$ws.Range("N140").Value = -95637.75

$ws.Range("H141").Value = 3199.8708
$ws.Range("I141").Value = 2056.7144
$ws.Range("K141").Value = 6170.1432
$ws.Range("M141").Value = -990.1431999999995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1596.6666
$ws.Range("I22").Value = 895
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 895
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -596
$ws.Range("N22").Value = -3598

$ws.Range("H32").Value = 15525.694
$ws.Range("I32").Value = 15736.461
$ws.Range("J32").Value = 13745.889
$ws.Range("K32").Value = 15736.461
$ws.Range("L32").Value = 13745.889
$ws.Range("M32").Value = -15449.461
$ws.Range("N32").Value = -14319.889

$ws.Range("H41").Value = 7566.6665
$ws.Range("I41").Value = 1450
$ws.Range("K41").Value = 1450
$ws.Range("M41").Value = -1036

$ws.Range("H45").Value = 1894.1333
$ws.Range("I45").Value = 1866.6666
$ws.Range("K45").Value = 1866.6666
$ws.Range("M45").Value = -1489.6666

$ws.Range("H61").Value = 5026.8237
$ws.Range("I61").Value = 3409.2432
$ws.Range("J61").Value = 9301.857
$ws.Range("K61").Value = 3409.2432
$ws.Range("L61").Value = 9301.857
$ws.Range("M61").Value = -3197.2432
$ws.Range("N61").Value = -9725.857

$ws.Range("H74").Value = 5658.393
$ws.Range("I74").Value = 2503.8235
$ws.Range("J74").Value = 10533.637
$ws.Range("K74").Value = 2503.8235
$ws.Range("L74").Value = 10533.637
$ws.Range("M74").Value = -1629.8235
$ws.Range("N74").Value = -12281.637

$ws.Range("H77").Value = 5658.393
$ws.Range("I77").Value = 2503.8235
$ws.Range("J77").Value = 10533.637
$ws.Range("K77").Value = 12519.1175
$ws.Range("L77").Value = 52668.185
$ws.Range("M77").Value = -8151.1175
$ws.Range("N77").Value = -61404.185

$ws.Range("H97").Value = 1140.8334
$ws.Range("I97").Value = 989
$ws.Range("J97").Value = 1900
$ws.Range("K97").Value = 989
$ws.Range("L97").Value = 1900
$ws.Range("M97").Value = -493
$ws.Range("N97").Value = -2892

$ws.Range("H132").Value = 3714.1355
$ws.Range("I132").Value = 1278.6342
$ws.Range("J132").Value = 9261.666999999999
$ws.Range("K132").Value = 3835.9026
$ws.Range("L132").Value = 27785.001
$ws.Range("M132").Value = -1305.9026
$ws.Range("N132").Value = -32845.001

$ws.Range("H136").Value = 5026.8237
$ws.Range("I136").Value = 3409.2432
$ws.Range("J136").Value = 9301.857
$ws.Range("K136").Value = 10227.7296
$ws.Range("L136").Value = 27905.571
$ws.Range("M136").Value = -7677.729599999999
$ws.Range("N136").Value = -33005.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 29800
$ws.Range("J33").Value = 29800
$ws.Range("L33").Value = 29800
$ws.Range("N33").Value = -30472

$ws.Range("H134").Value = 2843.3447
$ws.Range("I134").Value = 2917.842
$ws.Range("J134").Value = 2701.8
$ws.Range("K134").Value = 8753.526
$ws.Range("L134").Value = 8105.400000000001
$ws.Range("M134").Value = -6218.526
$ws.Range("N134").Value = -13175.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7140.577
$ws.Range("I31").Value = 13535.667
$ws.Range("K31").Value = 13535.667
$ws.Range("M31").Value = -13240.667

$ws.Range("H32").Value = 8970
$ws.Range("I32").Value = 8970
$ws.Range("K32").Value = 8970
$ws.Range("M32").Value = -8654

$ws.Range("H34").Value = 7140.577
$ws.Range("I34").Value = 13535.667
$ws.Range("K34").Value = 13535.667
$ws.Range("M34").Value = -13333.667

$ws.Range("H122").Value = 7811
$ws.Range("I122").Value = 4640.1665
$ws.Range("J122").Value = 14152.667
$ws.Range("K122").Value = 13920.4995
$ws.Range("L122").Value = 42458.001
$ws.Range("M122").Value = -11470.4995
$ws.Range("N122").Value = -47358.001

$ws.Range("H132").Value = 2051.4888
$ws.Range("I132").Value = 1500.6061
$ws.Range("J132").Value = 3566.4167
$ws.Range("K132").Value = 4501.8183
$ws.Range("L132").Value = 10699.2501
$ws.Range("M132").Value = -1971.8183
$ws.Range("N132").Value = -15759.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 9508.083000000001
$ws.Range("I51").Value = 540
$ws.Range("J51").Value = 15913.857
$ws.Range("K51").Value = 1620
$ws.Range("L51").Value = 47741.571
$ws.Range("M51").Value = -1160
$ws.Range("N51").Value = -48661.571

$ws.Range("H131").Value = 38080.117
$ws.Range("I131").Value = 2374.75
$ws.Range("J131").Value = 53949.168
$ws.Range("K131").Value = 7124.25
$ws.Range("L131").Value = 161847.504
$ws.Range("M131").Value = -2084.25
$ws.Range("N131").Value = -171927.504

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 7873.2
$ws.Range("I31").Value = 1110.3334
$ws.Range("K31").Value = 1110.3334
$ws.Range("M31").Value = -818.3334

$ws.Range("H37").Value = 7873.2
$ws.Range("I37").Value = 1110.3334
$ws.Range("K37").Value = 1110.3334
$ws.Range("M37").Value = -833.3334

$ws.Range("H46").Value = 29900
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 29900
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 29900
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -30212

$ws.Range("H51").Value = 21666.666
$ws.Range("J51").Value = 21666.666
$ws.Range("L51").Value = 21666.666
$ws.Range("N51").Value = -22684.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2179.25
$ws.Range("I16").Value = 2001
$ws.Range("J16").Value = 2238.6667
$ws.Range("K16").Value = 2001
$ws.Range("L16").Value = 2238.6667
$ws.Range("M16").Value = -1831
$ws.Range("N16").Value = -2578.6667

$ws.Range("H21").Value = 1699.2222
$ws.Range("I21").Value = 1699.2222
$ws.Range("K21").Value = 1699.2222
$ws.Range("M21").Value = -1525.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 70963.75
$ws.Range("J124").Value = 70963.75
$ws.Range("L124").Value = 70963.75
$ws.Range("N124").Value = -80783.75

$ws.Range("H132").Value = 1687.4166
$ws.Range("I132").Value = 847.58826
$ws.Range("J132").Value = 3727
$ws.Range("K132").Value = 2542.76478
$ws.Range("L132").Value = 11181
$ws.Range("M132").Value = -12.76477999999997
$ws.Range("N132").Value = -16241

$ws.Range("H136").Value = 5649.263
$ws.Range("I136").Value = 3668.0244
$ws.Range("K136").Value = 11004.0732
$ws.Range("M136").Value = -8454.073199999999
